$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds literal text (uses "." as a thousands
# separator, e.g. "24.580.50", and some values carry a significant
# trailing zero, e.g. "53.20"). Mark every Price cell we touch as Text
# before writing so Excel stores the literal string instead of
# re-parsing it as a number.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "24.580.50"
$ws.Range("E2").Value = "  +3.58%  "
$ws.Range("D3").Value = "1.695.11"
$ws.Range("E3").Value = "  +1.95%  "
$ws.Range("D5").Value = "316.44"
$ws.Range("E5").Value = "  +2.22%  "
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("E7").Value = "  +1.39%  "
$ws.Range("D8").Value = "0.4015"
$ws.Range("E8").Value = "  +1.93%  "
$ws.Range("D9").Value = "1.522"
$ws.Range("E10").Value = "  +0.06%  "
$ws.Range("D11").Value = "53.20"
$ws.Range("E11").Value = "  +7.16%  "
$ws.Range("E12").Value = "  +1.32%  "
$ws.Range("D13").Value = "7.224"
$ws.Range("E13").Value = "  +7.47%  "
$ws.Range("E14").Value = "  +2.85%  "
$ws.Range("D15").Value = "0.00001321"
$ws.Range("E15").Value = "  +0.70%  "
$ws.Range("D16").Value = "7.589"
$ws.Range("E16").Value = "  +5.19%  "
$ws.Range("D17").Value = "1.696.50"
$ws.Range("E17").Value = "  +1.77%  "
$ws.Range("D18").Value = "99.98"
$ws.Range("E18").Value = "  +0.44%  "
$ws.Range("D19").Value = "0.07054"
$ws.Range("E19").Value = "  +4.03%  "
$ws.Range("D20").Value = "19.67"
$ws.Range("E20").Value = "  +3.57%  "
$ws.Range("D21").Value = "6.875"
$ws.Range("E21").Value = "  +3.93%  "
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").Value = "14.06"
$ws.Range("E23").Value = "  +1.96%  "
$ws.Range("D24").Value = "24.579.74"
$ws.Range("E24").Value = "  +3.58%  "
$ws.Range("D25").Value = "3.033"
$ws.Range("E25").Value = "  +8.71%  "
$ws.Range("D26").Value = "2.320"
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("D27").Value = "22.41"
$ws.Range("E27").Value = "  +3.34%  "
$ws.Range("D28").Value = "160.07"
$ws.Range("E28").Value = "  +0.80%  "
$ws.Range("D29").Value = "5.220"
$ws.Range("E29").Value = "  +1.21%  "
$ws.Range("D30").Value = "134.60"
$ws.Range("E30").Value = "  +3.97%  "
$ws.Range("D31").Value = "7.441"
$ws.Range("E31").Value = "  +14.25%  "
$ws.Range("D32").Value = "1.883.38"
$ws.Range("E32").Value = "  +1.76%  "
$ws.Range("D33").Value = "1.104"
$ws.Range("E33").Value = "  -0.90%  "
$ws.Range("D34").Value = "0.08515"
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("D35").Value = "7.230"
$ws.Range("E35").Value = "  +9.24%  "
$ws.Range("D36").Value = "11.50"
$ws.Range("E36").Value = "  +11.00%  "
$ws.Range("D37").Value = "1.953"
$ws.Range("E37").Value = "  +0.23%  "
$ws.Range("D38").Value = "0.2729"
$ws.Range("E38").Value = "  +2.84%  "
$ws.Range("D39").Value = "14.49"
$ws.Range("E39").Value = "  +0.59%  "
$ws.Range("E40").Value = "  +9.26%  "
$ws.Range("D41").Value = "0.09067"
$ws.Range("E41").Value = "  +3.38%  "
$ws.Range("D42").Value = "1.462"
$ws.Range("E42").Value = "  +1.35%  "
$ws.Range("D43").Value = "0.7709"
$ws.Range("E43").Value = "  +2.67%  "
$ws.Range("D44").Value = "0.7214"
$ws.Range("E44").Value = "  +3.23%  "
$ws.Range("D45").Value = "2.543"
$ws.Range("E45").Value = "  +5.98%  "
$ws.Range("D46").Value = "15.42"
$ws.Range("E46").Value = "  +4.54%  "
$ws.Range("D47").Value = "4.210"
$ws.Range("E47").Value = "  +2.64%  "
$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").Value = "1.001"
$ws.Range("E48").Value = "  +0.11%  "
$ws.Range("B49").Value = "Flow"
$ws.Range("C49").Value = "https://coinranking.com/coin/QQ0NCmjVq+flow-flow"
$ws.Range("D49").Value = "1.345"
$ws.Range("E49").Value = "  +12.46%  "
$ws.Range("D50").Value = "141.33"
$ws.Range("E50").Value = "  +2.18%  "
$ws.Range("D51").Value = "0.08011"
$ws.Range("E51").Value = "  +3.31%  "
